$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 3 was "AP" ammo_338_federal; commit turns it into the "DMG" (hunting/
# Deer-hunter) loading of the same round, with updated ballistic data.

# Highlight the changed ammo type (A3:B3) and the now-notable damage/Rouble
# figure (E3) using the Accent2 theme color, matching the author's edit.
$ws.Range("A3:B3").Font.ThemeColor = 6
$ws.Range("E3").Font.ThemeColor = 6

# Update the ammo type label from AP to DMG.
$ws.Range("B3").Value = "DMG"

# Update price (Rouble) and damage (Joules) figures for the new loading.
$ws.Range("C3").Value = 4000
$ws.Range("H3").Value = 10.7

# Restore the cursor/selection position recorded in the saved workbook.
[void]$ws.Range("E11").Select()
